$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 74 values (B,C,D,F,G,H,I,J) - A74/E74 unchanged
$ws.Range("B74").Value = 613
$ws.Range("C74").Value = 232
$ws.Range("D74").Value = 382
$ws.Range("F74").Value = 184
$ws.Range("G74").Value = 130
$ws.Range("H74").Value = 497
$ws.Range("I74").Value = 90
$ws.Range("J74").Value = 1538

# Add new row 75 with the new quarter "01-04-2021" plus its figures.
# Force text format first so Excel doesn't auto-convert the date-like
# string into a date serial number, then restore the cell style so no
# extra explicit style is left on the cell.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "01-04-2021"
$ws.Range("A75").Style = "Normal"

$ws.Range("B75").Value = 655
$ws.Range("C75").Value = 226
$ws.Range("D75").Value = 429
$ws.Range("E75").Value = 1
$ws.Range("F75").Value = 182
$ws.Range("G75").Value = 128
$ws.Range("H75").Value = 421
$ws.Range("I75").Value = 99
$ws.Range("J75").Value = 1487
